$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").EntireColumn.Insert()

Write-Host "done"
